# Apply the "teacher free time hash calculated" edit described by the diff:
#  - Teachers sheet: update the course codes assigned to teachers AI and MMR
#  - TeacherFreeSlot sheet: update Monday/Wednesday/Thursday free-slot data
#  - Refresh the selection / active-cell bookmarks left behind on each sheet

$wb = $excel.ActiveWorkbook

$wsCourses         = $wb.Worksheets.Item("Courses")
$wsTeachers        = $wb.Worksheets.Item("Teachers")
$wsTeacherFreeSlot = $wb.Worksheets.Item("TeacherFreeSlot")

# --- Teachers sheet: swap in the course codes that replace "CSE 3203" / "CSE 1201" ---
$wsTeachers.Range("B2").Value = "CSE-1101"
$wsTeachers.Range("B3").Value = "EEE-1103"

# --- TeacherFreeSlot sheet: teacher AI is not free Monday, and the Wednesday/Thursday slot moved ---
$wsTeacherFreeSlot.Range("C2").Value = "no"
$wsTeacherFreeSlot.Range("E2").Value = "10.00-11.30;14.00-17.00"
$wsTeacherFreeSlot.Range("F2").Value = "10.00-11.30;14.00-17.00"

# --- Restore the selections left on each sheet (order matters: last one selected stays active) ---
[void]$wsCourses.Range("B4").Select()
[void]$wsTeachers.Range("C6").Select()
[void]$wsTeacherFreeSlot.Range("F2").Select()
